$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151; this pushes the existing rows
# 151-208 down to 152-209, preserving all of their original data.
$ws.Rows("151:151").Insert()

# Populate the newly inserted row 151 with a new weekly record.
$ws.Range("A151").Value = 8
$ws.Range("B151").Value = "Terminal La Palmera de La Serena"
$ws.Range("C151").Value = "Coquimbo"
$ws.Range("D151").Value = 44985
$ws.Range("D151").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E151").Value = 4
$ws.Range("F151").Value = 100112040
$ws.Range("G151").Value = "Cilantro"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 2400
$ws.Range("K151").Value = 2000
$ws.Range("L151").Value = 2500
$ws.Range("M151").Value = 2250
$ws.Range("N151").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O151").Value = "Provincia del Elquí"
$ws.Range("P151").Value = 1500
$ws.Range("Q151").Value = 1.5
$ws.Range("R151").Value = "Hortaliza"
